$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells (dots -> spaces, "Notes" -> "Notes_t")
$ws.Range("G1").Value = "Game 1"
$ws.Range("H1").Value = "Game 2"
$ws.Range("I1").Value = "Game 3"
$ws.Range("J1").Value = "Game 4"
$ws.Range("K1").Value = "Game 5"
$ws.Range("L1").Value = "Game 6"
$ws.Range("M1").Value = "Good Session"
$ws.Range("N1").Value = "Notes_t"

# Normalize lower-case "yes" entries to capitalized "Yes" (matches the existing "Yes" string used elsewhere)
$ws.Range("M3").Value = "Yes"
$ws.Range("M4").Value = "Yes"

# Select the header row in full (mirrors clicking the row-1 header) before turning on the
# AutoFilter, which is how Excel ends up leaving a full-row selection (A1:XFD1) behind.
$null = $ws.Rows("1:1").Select()

# Apply an AutoFilter across the data range
$null = $ws.Range("A1:N41").AutoFilter()

# Register the hidden _FilterDatabase defined name that Excel normally creates automatically
# alongside an AutoFilter.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$N`$41")
$filterName.Visible = $false
